# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap the "Islas Malvinas" / "Montserrat" rows (214 stays Groenlandia,
#     217 stays Santa Sede; rows 215/216 exchange labels + their D/H counts) ---
$ws.Range("A215").Value = "Montserrat"
$ws.Range("A216").Value = "Islas Malvinas"
$ws.Range("D215").Value = 12
$ws.Range("H215").Value = 1
$ws.Range("D216").Value = 13
$ws.Range("H216").Value = 0

# --- Refresh the "last updated" timestamp banner ---
$ws.Range("A1").Value = "Datos actualizados a 5 de Octubre de 2020 a las 14:10"

# --- Updated per-country case counts ---
# Iran (row 16)
$ws.Range("B16").Value = 475674
$ws.Range("C16").Value = 3902
$ws.Range("D16").Value = 392293
$ws.Range("E16").Value = 56189
$ws.Range("G16").Value = 235
$ws.Range("H16").Value = 27192

# Catar (row 36)
$ws.Range("B36").Value = 126692
$ws.Range("C36").Value = 194
$ws.Range("D36").Value = 123664
$ws.Range("E36").Value = 2812

# Kuwait (row 40)
$ws.Range("B40").Value = 107592
$ws.Range("C40").Value = 567
$ws.Range("D40").Value = 99549
$ws.Range("E40").Value = 7415
$ws.Range("G40").Value = 4
$ws.Range("H40").Value = 628

# Azerbaiyan (row 71)
$ws.Range("B71").Value = 40788
$ws.Range("C71").Value = 97
$ws.Range("D71").Value = 38587
$ws.Range("E71").Value = 1603
$ws.Range("G71").Value = 2
$ws.Range("H71").Value = 598

# Libia (row 75)
$ws.Range("B75").Value = 37437
$ws.Range("C75").Value = 628
$ws.Range("D75").Value = 22076
$ws.Range("E75").Value = 14765
$ws.Range("G75").Value = 4
$ws.Range("H75").Value = 596

# El Salvador (row 79)
$ws.Range("B79").Value = 29539
$ws.Range("C79").Value = 89
$ws.Range("D79").Value = 24406
$ws.Range("E79").Value = 4268

# Bosnia y Herzegovina (row 80)
$ws.Range("B80").Value = 28449
$ws.Range("C80").Value = 95
$ws.Range("D80").Value = 22032
$ws.Range("E80").Value = 5529
$ws.Range("G80").Value = 9
$ws.Range("H80").Value = 888

# Finlandia (row 104)
$ws.Range("B104").Value = 10702
$ws.Range("C104").Value = 164
$ws.Range("E104").Value = 2256
$ws.Range("G104").Value = 1
$ws.Range("H104").Value = 346

# Georgia (row 113)
$ws.Range("E113").Value = 4023
$ws.Range("G113").Value = 4
$ws.Range("H113").Value = 54

# Vietnam (row 168)
$ws.Range("B168").Value = 1097
$ws.Range("C168").Value = 1
$ws.Range("D168").Value = 1022
$ws.Range("E168").Value = 40
